$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5903
$ws.Range("K3").Value = 6088
$ws.Range("K4").Value = 1268
$ws.Range("K5").Value = 432
$ws.Range("K6").Value = 6688
$ws.Range("K7").Value = 20379

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 69
$ws.Range("K6").Value = 113
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 374
$ws.Range("K3").Value = 413
$ws.Range("K7").Value = 1347

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 324
$ws.Range("K4").Value = 43
$ws.Range("K6").Value = 260
$ws.Range("K7").Value = 881

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 117
$ws.Range("K3").Value = 121
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 344

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 200
$ws.Range("K3").Value = 228
$ws.Range("K5").Value = 30
$ws.Range("K7").Value = 691

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 155
$ws.Range("K7").Value = 475

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 138
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 180
$ws.Range("K7").Value = 593
$ws.Range("K8").Value = 1347
$ws.Range("K15").Value = 209
$ws.Range("K16").Value = 56
$ws.Range("K19").Value = 586
$ws.Range("K20").Value = 479
$ws.Range("K23").Value = 209
$ws.Range("K27").Value = 190
$ws.Range("K29").Value = 1110
$ws.Range("K30").Value = 76
$ws.Range("K31").Value = 228
$ws.Range("K33").Value = 881
$ws.Range("K34").Value = 114
$ws.Range("K36").Value = 259
$ws.Range("K37").Value = 691
$ws.Range("K42").Value = 758
$ws.Range("K44").Value = 172
$ws.Range("K49").Value = 111
$ws.Range("K51").Value = 260
$ws.Range("K52").Value = 538
$ws.Range("K53").Value = 261
$ws.Range("K54").Value = 395
$ws.Range("K55").Value = 226
$ws.Range("K63").Value = 60
$ws.Range("K65").Value = 475
$ws.Range("K67").Value = 800
$ws.Range("K68").Value = 55
$ws.Range("K73").Value = 183
$ws.Range("K74").Value = 23
$ws.Range("K77").Value = 144
$ws.Range("K78").Value = 229
$ws.Range("K79").Value = 505
$ws.Range("K84").Value = 162
$ws.Range("K85").Value = 950
$ws.Range("K88").Value = 220
$ws.Range("K89").Value = 297
$ws.Range("K90").Value = 188
$ws.Range("K94").Value = 273
$ws.Range("K95").Value = 344
$ws.Range("K97").Value = 162
$ws.Range("K98").Value = 99
$ws.Range("K99").Value = 335
$ws.Range("K101").Value = 20379

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 224
$ws.Range("K4").Value = 45
$ws.Range("K6").Value = 229
$ws.Range("K7").Value = 800

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 213
$ws.Range("K7").Value = 395

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 401
$ws.Range("K6").Value = 314
$ws.Range("K7").Value = 1110

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 176
$ws.Range("K7").Value = 586

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 206
$ws.Range("K4").Value = 31
$ws.Range("K6").Value = 282
$ws.Range("K7").Value = 758

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 70
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 226

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 111
$ws.Range("K6").Value = 49

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 165
$ws.Range("K7").Value = 505

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 160
$ws.Range("K7").Value = 479

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 76
$ws.Range("K7").Value = 259

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 193
$ws.Range("K4").Value = 22
$ws.Range("K6").Value = 159
$ws.Range("K7").Value = 593

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 273

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 45
$ws.Range("K4").Value = 14
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 297

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 52
$ws.Range("K7").Value = 190

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 69
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 260

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 326
$ws.Range("K7").Value = 950

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 145
$ws.Range("K6").Value = 189
$ws.Range("K7").Value = 538

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 23
